# Aula14/tabelinha.xlsx — add three new MUX columns to the truth table and
# refresh the data/view to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three new columns at the positions where they appear in the
#    final layout (F, H, J). Each Insert() shifts everything at/after the
#    target column one slot to the right, which is exactly how this was
#    built by hand in Excel (habEscritaReg/Operacao/BEQ/... all slide over).
# ---------------------------------------------------------------------
$ws.Columns("F").Insert()
$ws.Columns("H").Insert()
$ws.Columns("J").Insert()

# ---------------------------------------------------------------------
# 2. Header row (row 2) text for the three new columns.
#    Insertion order controls shared-string table order, matching the
#    target file (13=MUX ULA mem, 14=MUX Rt/Imediado, 15=MUX Rt/Rd).
# ---------------------------------------------------------------------
$ws.Range("J2").Value = "MUX ULA mem"
$ws.Range("H2").Value = "MUX Rt/Imediado"
$ws.Range("F2").Value = "MUX Rt/Rd"

# ---------------------------------------------------------------------
# 3. Data values for the new columns (rows 3-5) plus the updated values
#    in the columns that kept their meaning but changed numbers.
# ---------------------------------------------------------------------
# Row 3 (LW)
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 1

# Row 4 (SW)
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = 0

# Row 5 (BEQ)
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1

# ---------------------------------------------------------------------
# 4. Column widths for the three freshly-inserted columns (they come in
#    with no explicit width / default width).
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 10.67
$ws.Columns("H").ColumnWidth = 15
$ws.Columns("J").ColumnWidth = 14.67

# ---------------------------------------------------------------------
# 5. Sheet view: zoom out a bit and move the visible window / selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 80
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("F3:M3").Select()
